# Harmonize similar tags to be the same.
# On the "isa_template" sheet, the Tags entry (row 12-14) previously stored
# "extraction protocol" / "MINSEQ" / "Extraction" with accession "DPBO:1000171"
# and source ref "DPBO". Update it to use the harmonized tag
# "extraction" / "RNA" with accession "OBI:0302884" (no term source ref).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("isa_template")

$ws.Range("C12").Value = "RNA"
$ws.Range("B12").Value = "extraction"
$ws.Range("D12").ClearContents()

$ws.Range("B13").Value = "OBI:0302884"

$ws.Range("B14").ClearContents()

$ws.Activate()
$ws.Range("D13").Select()
